# Weekly update: two new price-report rows for "Ají" (Americana) are added to the dataset.
# The source table is appended chronologically as new market-report rows arrive each week,
# so inserting at a specific row (rather than appending at the very end) matches how the
# published consolidated file grows over time.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert first new row at position 31 and populate it.
$ws.Rows("31:31").Insert()
$ws.Range("A31").Value = 7
$ws.Range("B31").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C31").Value = 'Ñuble'
$ws.Range("D31").Value = 44664
$ws.Range("E31").Value = 16
$ws.Range("F31").Value = 100112021
$ws.Range("G31").Value = 'Ají'
$ws.Range("H31").Value = 'Americana (o)'
$ws.Range("I31").Value = 'Primera'
$ws.Range("J31").Value = 80
$ws.Range("K31").Value = 8500
$ws.Range("L31").Value = 9000
$ws.Range("M31").Value = 8750
$ws.Range("N31").Value = '$/caja 15 kilos'
$ws.Range("O31").Value = 'Región del Maule'
$ws.Range("P31").Value = 583
$ws.Range("Q31").Value = 15
$ws.Range("R31").Value = 'Hortaliza'

# Insert second new row; after the first insert shifted everything down, this new row lands at 64.
$ws.Rows("64:64").Insert()
$ws.Range("A64").Value = 7
$ws.Range("B64").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C64").Value = 'Ñuble'
$ws.Range("D64").Value = 44663
$ws.Range("E64").Value = 16
$ws.Range("F64").Value = 100112021
$ws.Range("G64").Value = 'Ají'
$ws.Range("H64").Value = 'Americana (o)'
$ws.Range("I64").Value = 'Primera'
$ws.Range("J64").Value = 80
$ws.Range("K64").Value = 8500
$ws.Range("L64").Value = 9000
$ws.Range("M64").Value = 8750
$ws.Range("N64").Value = '$/caja 15 kilos'
$ws.Range("O64").Value = 'Región del Maule'
$ws.Range("P64").Value = 583
$ws.Range("Q64").Value = 15
$ws.Range("R64").Value = 'Hortaliza'
